$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (row 42) down onto the
# two new rows (43, 44) so the new cells pick up the same styles (bold
# centered bordered index column, date-time number format on the match-date
# column, etc.) without hard-coding style ids.
$ws.Range("A42:V42").Copy()
$ws.Range("A43:V44").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Row 43 ----
$ws.Cells.Item(43, 1).Value = 42
$ws.Cells.Item(43, 2).Value = "armenia"
$ws.Cells.Item(43, 3).Value = "premier-league"
$ws.Cells.Item(43, 4).Value = "2023-2024"
$ws.Cells.Item(43, 5).Value = 45194.58333333334
$ws.Cells.Item(43, 6).Value = "Pyunik Yerevan"
$ws.Cells.Item(43, 7).Value = 1
$ws.Cells.Item(43, 8).Value = "Shirak Gyumri"
$ws.Cells.Item(43, 9).Value = 0
$ws.Cells.Item(43, 10).Value = 1.12
$ws.Cells.Item(43, 11).Value = "24/09/2023 01:12"
$ws.Cells.Item(43, 12).Value = 1.12
$ws.Cells.Item(43, 13).Value = "25/09/2023 13:52"
$ws.Cells.Item(43, 14).Value = 7.86
$ws.Cells.Item(43, 15).Value = "24/09/2023 01:12"
$ws.Cells.Item(43, 16).Value = 8.98
$ws.Cells.Item(43, 17).Value = "25/09/2023 13:52"
$ws.Cells.Item(43, 18).Value = 14.85
$ws.Cells.Item(43, 19).Value = "24/09/2023 01:12"
$ws.Cells.Item(43, 20).Value = 19.72
$ws.Cells.Item(43, 21).Value = "25/09/2023 13:53"
$ws.Cells.Item(43, 22).Value = "https://www.betexplorer.com/football/armenia/premier-league/pyunik-yerevan-shirak-gyumri/GlILN5hj/"

# ---- Row 44 ----
$ws.Cells.Item(44, 1).Value = 43
$ws.Cells.Item(44, 2).Value = "armenia"
$ws.Cells.Item(44, 3).Value = "premier-league"
$ws.Cells.Item(44, 4).Value = "2023-2024"
$ws.Cells.Item(44, 5).Value = 45194.66666666666
$ws.Cells.Item(44, 6).Value = "Ararat-Armenia"
$ws.Cells.Item(44, 7).Value = 3
$ws.Cells.Item(44, 8).Value = "Alashkert"
$ws.Cells.Item(44, 9).Value = 1
$ws.Cells.Item(44, 10).Value = 2.31
$ws.Cells.Item(44, 11).Value = "24/09/2023 03:13"
$ws.Cells.Item(44, 12).Value = 2.1
$ws.Cells.Item(44, 13).Value = "25/09/2023 15:59"
$ws.Cells.Item(44, 14).Value = 3.25
$ws.Cells.Item(44, 15).Value = "24/09/2023 03:13"
$ws.Cells.Item(44, 16).Value = 3.32
$ws.Cells.Item(44, 17).Value = "25/09/2023 15:59"
$ws.Cells.Item(44, 18).Value = 2.85
$ws.Cells.Item(44, 19).Value = "24/09/2023 03:13"
$ws.Cells.Item(44, 20).Value = 3.57
$ws.Cells.Item(44, 21).Value = "25/09/2023 15:59"
$ws.Cells.Item(44, 22).Value = "https://www.betexplorer.com/football/armenia/premier-league/ararat-armenia-alashkert/bDNGOowp/"
